# chore: adapt column header formatting to respective input file names
#
# - Rename the "_old"/"_new" header-column suffixes to the respective
#   format-version identifiers ("_FV2310" / "_FV2404").
# - Turn the sheet's data range into an actual Excel Table (ListObject)
#   so the new headers get an AutoFilter + table formatting.
# - Freeze the header row so it stays visible while scrolling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header cells in row 1 -------------------------------------
# Columns A:J were suffixed "_old" (now "_FV2310"), columns L:U were
# suffixed "_new" (now "_FV2404"); column K ("diff") is untouched.
$lastCol = $ws.UsedRange.Columns.Count
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $header = [string]$cell.Value2
    if ($header.EndsWith("_old")) {
        $cell.Value = $header.Substring(0, $header.Length - 4) + "_FV2310"
    } elseif ($header.EndsWith("_new")) {
        $cell.Value = $header.Substring(0, $header.Length - 4) + "_FV2404"
    }
}

# --- 2. Convert the used range into an Excel Table ------------------------
$lastRow = $ws.UsedRange.Rows.Count
$tableRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# --- 3. Freeze the header row ---------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
